$wb = $excel.ActiveWorkbook

# The credentials used by "script 1" live on the Login sheet.
$wb.Worksheets.Item("Login").Activate()
$ws = $wb.ActiveSheet

# Swap the sample login from the old "ajenkins" account over to the
# AGSAuto service account (commit: "change script 1 data to run on
# AGSAuto account instead of ajenkins").
$ws.Range("A2").Value = "AGSAutoT02"
$ws.Range("B2").Value = "SERVICE`$08"

# Re-fit the columns now that the new values have different lengths.
$ws.Range("A1:F2").EntireColumn.AutoFit() | Out-Null

# Leave the cursor where the user clicked next.
$ws.Range("B3").Select()
